# Atualiza datasets e ajustes das ligas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Rodada" headers (B1:T1) -> Rodada 1 .. Rodada 19
$headers = @(
    "Rodada 1","Rodada 2","Rodada 3","Rodada 4","Rodada 5","Rodada 6","Rodada 7",
    "Rodada 8","Rodada 9","Rodada 10","Rodada 11","Rodada 12","Rodada 13","Rodada 14",
    "Rodada 15","Rodada 16","Rodada 17","Rodada 18","Rodada 19"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i  # column B = 2
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Full alphabetically sorted team list (14 existing + 4 new teams) for rows 2..19
$teams = @(
    "bugredasmissões",
    "C R Juvenal",
    "Doug Leal F.C",
    "Esquadrão Gazembrino",
    "FBC Colorado",
    "GaúchoDaFronteira F.C",
    "GE Bebum",
    "GrioTeam",
    "Grêmio_Campeão_LA_27",
    "JV5 Tricolor Gaúcho",
    "La Primeira Patada Es Nuestra",
    "lsauer fc",
    "Medonho´s F.C.",
    "NHU PORÃ SAF.",
    "Pontaç0 F.C.",
    "SC 100 Sono",
    "SC ÉoINTER!",
    "Texas Club 2026"
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $teams[$i]

    if ($row -gt 15) {
        # New rows (16..19): copy formatting from an existing team row and fill B:T with 0
        $srcRow = $ws.Range($ws.Cells.Item(15, 1), $ws.Cells.Item(15, 20))
        $dstRow = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 20))
        $srcRow.Copy() | Out-Null
        $dstRow.PasteSpecial(-4122) | Out-Null # xlPasteFormats

        $dataRange = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 20))
        $dataRange.Value = 0
    }
}

$ws.UsedRange | Out-Null
